# chore: simulator full-month coverage, persist logs, fix employees
$wb = $excel.ActiveWorkbook

# --- Fix employee/client names (typo fixes from the simulator) ---
$wsTime = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# Weekly Timesheet: client names in column B
$wsTime.Range("B4").Value = "Tubergen"
$wsTime.Range("B5").Value = "Hewett"

# Jason Schema: mirrors the same rows in column D
$wsSchema.Range("D4").Value = "Tubergen"
$wsSchema.Range("D5").Value = "Hewett"

# Employee ID fix (Jason Schema column B, all rows reference the same employee)
$wsSchema.Range("B2:B5").Value = "emp_35u1tnme"

# --- Full-month coverage: populate Rate/Total for the simulated rows ---
# Weekly Timesheet: E=Rate, F=Total
$wsTime.Range("E2").Value = 150
$wsTime.Range("F2").Value = 1200
$wsTime.Range("E3").Value = 150
$wsTime.Range("F3").Value = 1200
$wsTime.Range("E4").Value = 150
$wsTime.Range("F4").Value = 3000
$wsTime.Range("E5").Value = 150
$wsTime.Range("F5").Value = 3000

# Subtotal / Hourly subtotal / Grand total (Weekly Timesheet column F)
$wsTime.Range("F7").Value = 8400
$wsTime.Range("F11").Value = 8400
$wsTime.Range("F12").Value = 8400

# Jason Schema: F=Rate, G=Total
$wsSchema.Range("F2").Value = 150
$wsSchema.Range("G2").Value = 1200
$wsSchema.Range("F3").Value = 150
$wsSchema.Range("G3").Value = 1200
$wsSchema.Range("F4").Value = 150
$wsSchema.Range("G4").Value = 3000
$wsSchema.Range("F5").Value = 150
$wsSchema.Range("G5").Value = 3000
